$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cthrc1"
$ws.Range("C2").Value = "Fzd3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 5.507882
$ws.Range("H2").Value = 16.523646
$ws.Range("I2").Value = 0.03518866199235487
$ws.Range("J2").Value = 0.03518866199235487
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2590246666666667
$ws.Range("N2").Value = 0.777074
$ws.Range("O2").Value = 0.1073177818850196
$ws.Range("P2").Value = 0.1073177818850196
$ws.Range("Q2").Value = 1.426677299089333
$ws.Range("R2").Value = 12.840095691804
$ws.Range("S2").Value = 0.00377636915252122
$ws.Range("T2").Value = 0.00377636915252122

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Cthrc1"
$ws.Range("C3").Value = "Fzd3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 5.507882
$ws.Range("H3").Value = 16.523646
$ws.Range("I3").Value = 0.03518866199235487
$ws.Range("J3").Value = 0.03518866199235487
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.8886716666666666
$ws.Range("N3").Value = 2.666015
$ws.Range("O3").Value = 0.3681899230603399
$ws.Range("P3").Value = 0.3681899230603398
$ws.Range("Q3").Value = 4.894698676743332
$ws.Range("R3").Value = 44.05228809069
$ws.Range("S3").Value = 0.01295611075156144
$ws.Range("T3").Value = 0.01295611075156144

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Cthrc1"
$ws.Range("C4").Value = "Fzd3"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 5.507882
$ws.Range("H4").Value = 16.523646
$ws.Range("I4").Value = 0.03518866199235487
$ws.Range("J4").Value = 0.03518866199235487
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.000138
$ws.Range("N4").Value = 0.000414
$ws.Range("O4").Value = 0.0000571754578076195
$ws.Range("P4").Value = 0.0000571754578076195
$ws.Range("Q4").Value = 0.0007600877159999999
$ws.Range("R4").Value = 0.006840789443999999
$ws.Range("S4").Value = 0.00000201192785905047
$ws.Range("T4").Value = 0.00000201192785905047

$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Cthrc1"
$ws.Range("C5").Value = "Fzd3"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 5.507882
$ws.Range("H5").Value = 16.523646
$ws.Range("I5").Value = 0.03518866199235487
$ws.Range("J5").Value = 0.03518866199235487
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.265788666666667
$ws.Range("N5").Value = 3.797366
$ws.Range("O5").Value = 0.524435119596833
$ws.Range("P5").Value = 0.524435119596833
$ws.Range("Q5").Value = 6.971814612937333
$ws.Range("R5").Value = 62.746331516436
$ws.Range("S5").Value = 0.01845417016041316
$ws.Range("T5").Value = 0.01845417016041316

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Cthrc1"
$ws.Range("C6").Value = "Fzd3"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 149.9875183333334
$ws.Range("H6").Value = 449.9625550000001
$ws.Range("I6").Value = 0.9582376829612175
$ws.Range("J6").Value = 0.9582376829612176
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.2590246666666667
$ws.Range("N6").Value = 0.777074
$ws.Range("O6").Value = 0.1073177818850196
$ws.Range("P6").Value = 0.1073177818850196
$ws.Range("Q6").Value = 38.85046694045223
$ws.Range("R6").Value = 349.6542024640701
$ws.Range("S6").Value = 0.1028359426540385
$ws.Range("T6").Value = 0.1028359426540385

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Cthrc1"
$ws.Range("C7").Value = "Fzd3"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 149.9875183333334
$ws.Range("H7").Value = 449.9625550000001
$ws.Range("I7").Value = 0.9582376829612175
$ws.Range("J7").Value = 0.9582376829612176
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.8886716666666666
$ws.Range("N7").Value = 2.666015
$ws.Range("O7").Value = 0.3681899230603399
$ws.Range("P7").Value = 0.3681899230603398
$ws.Range("Q7").Value = 133.2896578964806
$ws.Range("R7").Value = 1199.606921068325
$ws.Range("S7").Value = 0.352813458763009
$ws.Range("T7").Value = 0.352813458763009

$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Cthrc1"
$ws.Range("C8").Value = "Fzd3"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 149.9875183333334
$ws.Range("H8").Value = 449.9625550000001
$ws.Range("I8").Value = 0.9582376829612175
$ws.Range("J8").Value = 0.9582376829612176
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.000138
$ws.Range("N8").Value = 0.000414
$ws.Range("O8").Value = 0.0000571754578076195
$ws.Range("P8").Value = 0.0000571754578076195
$ws.Range("Q8").Value = 0.02069827753
$ws.Range("R8").Value = 0.18628449777
$ws.Range("S8").Value = 0.00005478767821182016
$ws.Range("T8").Value = 0.00005478767821182017

$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Cthrc1"
$ws.Range("C9").Value = "Fzd3"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 149.9875183333334
$ws.Range("H9").Value = 449.9625550000001
$ws.Range("I9").Value = 0.9582376829612175
$ws.Range("J9").Value = 0.9582376829612176
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.265788666666667
$ws.Range("N9").Value = 3.797366
$ws.Range("O9").Value = 0.524435119596833
$ws.Range("P9").Value = 0.524435119596833
$ws.Range("Q9").Value = 189.8525008477922
$ws.Range("R9").Value = 1708.67250763013
$ws.Range("S9").Value = 0.5025334938659582
$ws.Range("T9").Value = 0.5025334938659582

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Cthrc1"
$ws.Range("C10").Value = "Fzd3"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.028937
$ws.Range("H10").Value = 3.086811
$ws.Range("I10").Value = 0.006573655046427582
$ws.Range("J10").Value = 0.006573655046427582
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.2590246666666667
$ws.Range("N10").Value = 0.777074
$ws.Range("O10").Value = 0.1073177818850196
$ws.Range("P10").Value = 0.1073177818850196
$ws.Range("Q10").Value = 0.266520063446
$ws.Range("R10").Value = 2.398680571014
$ws.Range("S10").Value = 0.0007054700784598737
$ws.Range("T10").Value = 0.0007054700784598737

$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Cthrc1"
$ws.Range("C11").Value = "Fzd3"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1.028937
$ws.Range("H11").Value = 3.086811
$ws.Range("I11").Value = 0.006573655046427582
$ws.Range("J11").Value = 0.006573655046427582
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.8886716666666666
$ws.Range("N11").Value = 2.666015
$ws.Range("O11").Value = 0.3681899230603399
$ws.Range("P11").Value = 0.3681899230603398
$ws.Range("Q11").Value = 0.914387158685
$ws.Range("R11").Value = 8.229484428165
$ws.Range("S11").Value = 0.002420353545769386
$ws.Range("T11").Value = 0.002420353545769386

$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Cthrc1"
$ws.Range("C12").Value = "Fzd3"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1.028937
$ws.Range("H12").Value = 3.086811
$ws.Range("I12").Value = 0.006573655046427582
$ws.Range("J12").Value = 0.006573655046427582
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.000138
$ws.Range("N12").Value = 0.000414
$ws.Range("O12").Value = 0.0000571754578076195
$ws.Range("P12").Value = 0.0000571754578076195
$ws.Range("Q12").Value = 0.000141993306
$ws.Range("R12").Value = 0.001277939754
$ws.Range("S12").Value = 0.0000003758517367488652
$ws.Range("T12").Value = 0.0000003758517367488652

$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Cthrc1"
$ws.Range("C13").Value = "Fzd3"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1.028937
$ws.Range("H13").Value = 3.086811
$ws.Range("I13").Value = 0.006573655046427582
$ws.Range("J13").Value = 0.006573655046427582
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 1.265788666666667
$ws.Range("N13").Value = 3.797366
$ws.Range("O13").Value = 0.524435119596833
$ws.Range("P13").Value = 0.524435119596833
$ws.Range("Q13").Value = 1.302416793314
$ws.Range("R13").Value = 11.721751139826
$ws.Range("S13").Value = 0.003447455570461573
$ws.Range("T13").Value = 0.003447455570461573
